$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Dinsdag 15 november 2012" -> "Dinsdag 13 november 2012"
# ---------------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("15 november 2012", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "13 november 2012", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Append the new diary entry (4 paragraphs) after the last paragraph,
#    moving the _GoBack bookmark into the new heading paragraph.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Remove the stray _GoBack bookmark sitting at the end of the last paragraph;
# it will be re-created further down, inside the new heading line.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rng = $lastPara.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()            # blank separator paragraph

$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()              # heading: "Maandag 19 november 2012"

$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()              # "8;45 - 12;00"

$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()              # body text paragraph

$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count - 2)
$headingPara.Range.Text = "Maandag 19 november 2012"
$headingPara.Range.Font.Bold = $true
$headingPara.Range.Font.Underline = 1

$timePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$timePara.Range.Text = "8;45 " + [char]0x2013 + " 12;00"
$timePara.Range.Font.Bold = $true
$timePara.Range.Font.Underline = 1

$bodyPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bodyPara.Range.Text = "Mail functie toegevoegd aan de website. De persoon die zich wilt registreren krijgt een activatie mail met daarin een link waarop hij kan klikken, zodat hij op de activatie.php pagina komt. Hier kan hij zijn MD5-hash wachtwoord veranderen naar een wachtwoord dat hij zelf wilt. We hebben twee soorten mail gemaakt, een met plain tekst en een met HTML text. "

# Re-insert the _GoBack bookmark inside the heading paragraph, between
# "Maandag 19 novembe" and "r 2012", as in the authored edit.
$headingStart = $headingPara.Range.Start
$bookmarkPos = $headingStart + ("Maandag 19 novembe").Length
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
